$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.211098670959473
$ws.Range("B1").Value = 3.1682288646698
$ws.Range("C1").Value = 2.663084030151367
$ws.Range("D1").Value = 2.222048044204712
$ws.Range("E1").Value = 1.496410250663757
